$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.437.90"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.573.93"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.31"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  +2.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.88"
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3426"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.161"
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07669"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.010"
$ws.Range("E14").Value = "  -0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.939"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "1.570.63"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001134"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.34"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06755"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.82"
$ws.Range("E21").Value = "  +2.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.228"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.05"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.428"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("B25").Value = "WrappedBTC"
$ws.Range("C25").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D25").Value = "22.428.97"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -8.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.34"
$ws.Range("E27").Value = "  +2.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "146.10"
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.033"
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.37"
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("D31").Value = "1.747.84"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.214"
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.017"
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.002"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.05"
$ws.Range("E35").Value = "  -2.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08582"
$ws.Range("E36").Value = "  +1.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02550"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2319"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06597"
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.340"
$ws.Range("E40").Value = "  +7.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.481"
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6458"
$ws.Range("E42").Value = "  +1.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.57"
$ws.Range("E43").Value = "  -1.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.20"
$ws.Range("E44").Value = "  -2.84%  "
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.800"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6018"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.299"
$ws.Range("E48").Value = "  +8.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.087"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "125.65"
$ws.Range("E50").Value = "  +3.71%  "
$ws.Range("E51").Value = "  +0.72%  "
